$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Add a new tracking entry to Table4 (M1:Q15 -> M1:Q16) ---------------
$lo4 = $ws.ListObjects.Item("Table4")

# Copy the formatting of the current last row (row 15) down onto the new
# row 16 so the new row looks like the rest of the tracker (borders,
# number formats, etc.) before we fill in the new values.
$ws.Range("M15:Q15").Copy()
$ws.Range("M16").PasteSpecial(-4122)

$ws.Range("M16").Value = "Problem Solving(Algorithms & Data Structures)"
$ws.Range("N16").Value = 45101
$ws.Range("O16").Value = "1067.97/2200"
$ws.Range("P16").Value = 143690
$ws.Range("Q16").Formula = "=IF(ROW()>2,(`$P`$2-P16)/`$P`$2,""NA"")"

# Grow the table definition to include the freshly-entered row.
$lo4.Resize($ws.Range("M1:Q16"))

# --- 2) Remove the blank spacer row that used to sit above the Python /  ---
# --- C++ category tables (row 19); everything below shifts up by one.   ---
$ws.Rows("19:19").Delete()

# --- 3) Restore the cursor position the author left the sheet in. ----------
$ws.Range("P16").Select()
